$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey_6")

# Insert a new row at position 78 (shifts existing rows 78+ down by one)
$ws.Rows.Item(78).Insert()

$ws.Cells.Item(78, 2).Value = "hhm_elevated_risk"
$ws.Cells.Item(78, 1).Value = "q28a_scale_1"
